# Applies the CLIP Model Results update described by the commit:
# "been working on clip_openai script for compute-canada"
#
# The edit appends 8 new experiment-result rows (rows 19-26) to Sheet1,
# describing new runs of the clip_openai script on compute-canada
# (slurm job results), and moves the active selection to reflect where
# the author was last working (I22), scrolled so row 17 is visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = "1-7 (training), 0 (testing)"
$ws.Range("B19").Value = "8 words"
$ws.Range("C19").Value = "3D fMRI image 4 seconds after last word"
$ws.Range("D19").Value = "Default CLIP"
$ws.Range("E19").Value = "Default CLIP"
$ws.Range("F19").Value = "3D Resnet18"
$ws.Range("G19").Value = "Cosine Similarity"
$ws.Range("H19").Value = "embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=24, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I19").Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = "700/0/100"
$ws.Range("L19").Value = 0.42389838129496399
$ws.Range("M19").Value = 0.439186151079136
$ws.Range("N19").Value = 0.11939102564102499
$ws.Range("O19").Value = 0.13541666666666599
$ws.Range("P19").Value = "slurm-42245994"

# Row 20
$ws.Range("A20").Value = "1-7 (training), 0 (testing)"
$ws.Range("B20").Value = "8 words"
$ws.Range("C20").Value = "3D fMRI image 4 seconds after last word"
$ws.Range("D20").Value = "Default CLIP"
$ws.Range("E20").Value = "Default CLIP"
$ws.Range("F20").Value = "3D Resnet18"
$ws.Range("G20").Value = "Cosine Similarity"
$ws.Range("H20").Value = "embed_dim=512, image_resolution, layers=(2,2,2,2), width=64, context_length=24, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I20").Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = "700/0/100"
$ws.Range("L20").Value = 0.42153776978417201
$ws.Range("M20").Value = 0.42884442446043097
$ws.Range("N20").Value = 0.0881410256410256
$ws.Range("O20").Value = 0.103365384615384
$ws.Range("P20").Value = "slurm-42257491"

# Row 21
$ws.Range("A21").Value = "1-7 (training), 0 (testing)"
$ws.Range("B21").Value = "8 words"
$ws.Range("C21").Value = "(Averaged over all subjects) 3D fMRI image 4 seconds after last word"
$ws.Range("D21").Value = "Default CLIP"
$ws.Range("E21").Value = "Default CLIP"
$ws.Range("F21").Value = "3D Resnet18"
$ws.Range("G21").Value = "Cosine Similarity"
$ws.Range("H21").Value = "embed_dim=512, image_resolution, layers=(2,2,2,2), width=64, context_length=24, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I21").Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = "100/0/100"
$ws.Range("L21").Value = 0.99278846153846101
$ws.Range("M21").Value = 0.98557692307692302
$ws.Range("N21").Value = 0.0657051282051282
$ws.Range("O21").Value = 0.0665064102564102
$ws.Range("P21").Value = "slurm-42261175"

# Row 22
$ws.Range("A22").Value = "1-7 (training), 0 (testing)"
$ws.Range("B22").Value = "8 words"
$ws.Range("C22").Value = "(Averaged over all subjects) 3D fMRI image 4 seconds after last word"
$ws.Range("D22").Value = "Default CLIP"
$ws.Range("E22").Value = "Default CLIP"
$ws.Range("F22").Value = "3D Resnet18"
$ws.Range("G22").Value = "Cosine Similarity"
$ws.Range("H22").Value = "embed_dim=512, image_resolution, layers=(2,2,2,2), width=64, context_length=24, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I22").Value = "LR=1e-4, batch_size=32, weight_decay=0.2"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = "100/0/100"
$ws.Range("P22").Value = "slurm-42263147"

# Row 23
$ws.Range("A23").Value = "1-7 (training), 0 (testing)"
$ws.Range("B23").Value = "8 words"
$ws.Range("C23").Value = "(Averaged over all subjects) 3D fMRI image 4 seconds after last word"
$ws.Range("D23").Value = "Default CLIP"
$ws.Range("E23").Value = "Default CLIP"
$ws.Range("F23").Value = "3D Resnet18"
$ws.Range("G23").Value = "Cosine Similarity"
$ws.Range("H23").Value = "embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=24, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I23").Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = "100/0/100"
$ws.Range("P23").Value = "slurm-42263522"

# Row 24
$ws.Range("A24").Value = "1-7 (training), 0 (testing)"
$ws.Range("B24").Value = "8 words"
$ws.Range("C24").Value = "(Averaged over all subjects) 3D fMRI image 4 seconds after last word"
$ws.Range("D24").Value = "Default CLIP"
$ws.Range("E24").Value = "Default CLIP"
$ws.Range("F24").Value = "3D Resnet18"
$ws.Range("G24").Value = "Cosine Similarity"
$ws.Range("H24").Value = "embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=24, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I24").Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Range("J24").Value = 10
$ws.Range("K24").Value = "100/0/100"
$ws.Range("P24").Value = "slurm-42263524"

# Row 25
$ws.Range("A25").Value = "1-7 (training), 0 (testing)"
$ws.Range("B25").Value = "16 words"
$ws.Range("C25").Value = "(Averaged over all subjects) 3D fMRI image 6 seconds after last word"
$ws.Range("D25").Value = "Default CLIP"
$ws.Range("E25").Value = "Default CLIP"
$ws.Range("F25").Value = "3D Resnet18"
$ws.Range("G25").Value = "Cosine Similarity"
$ws.Range("H25").Value = "embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=48, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I25").Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = "100/0/100"
$ws.Range("P25").Value = "slurm-42263542"

# Row 26
$ws.Range("A26").Value = "1-7 (training), 0 (testing)"
$ws.Range("B26").Value = "4 words"
$ws.Range("C26").Value = "(Averaged over all subjects) (detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word"
$ws.Range("D26").Value = "Default CLIP"
$ws.Range("E26").Value = "Default CLIP"
$ws.Range("F26").Value = "3D Resnet18"
$ws.Range("G26").Value = "Cosine Similarity"
$ws.Range("H26").Value = "embed_dim=512, image_resolution, layers=(2,2,2,2), width=64, context_length=24, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Range("I26").Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = "100/0/100"

# Reflect where the author was last working when the file was saved.
$ws.Range("I22").Select()
